$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("A2").Value = "Tổng công tại CẦN THƠ"
$ws.Range("B2").Value = 28
$ws.Range("A3").Value = "Phụ cấp tại CẦN THƠ"
$ws.Range("B3").Value = 980000
$ws.Range("A4").Value = "Lương cơ bản tại CẦN THƠ"
$ws.Range("B4").Value = 17250000
$ws.Range("A5").Value = "Chiết khấu sale chính tại CẦN THƠ"
$ws.Range("B5").Value = 1000000
$ws.Range("A6").Value = "Chiết khấu sale phụ tại CẦN THƠ"
$ws.Range("B6").Value = 0
$ws.Range("A7").Value = "Đơn 1 bác sĩ tại CẦN THƠ"
$ws.Range("B7").Value = 0
$ws.Range("A8").Value = "Đơn 2 bác sĩ tại CẦN THƠ"
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = "Công phụ phẫu 1 tại CẦN THƠ"
$ws.Range("B9").Value = 0
$ws.Range("A10").Value = "Công phụ phẫu 2 tại CẦN THƠ"
$ws.Range("B10").Value = 0
$ws.Range("A11").Value = "Chiết khấu thu nợ tại CẦN THƠ"
$ws.Range("B11").Value = 510000.0000000001
$ws.Range("A12").Value = "Ứng lương tại CẦN THƠ"
$ws.Range("B12").Value = -4469000
$ws.Range("A13").Value = "Tổng công tại LONG XUYÊN"
$ws.Range("B13").Value = 0
$ws.Range("A14").Value = "Phụ cấp tại LONG XUYÊN"
$ws.Range("B14").Value = 0
$ws.Range("A15").Value = "Lương công tác tại LONG XUYÊN"
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = "Lương cơ bản tại LONG XUYÊN"
$ws.Range("B16").Value = 11500000
$ws.Range("A17").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws.Range("B17").Value = 910000
$ws.Range("A18").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws.Range("B18").Value = 0
$ws.Range("A19").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws.Range("B19").Value = 0
$ws.Range("A20").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws.Range("B20").Value = 0
$ws.Range("A21").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws.Range("B21").Value = 0
$ws.Range("A22").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws.Range("B22").Value = 0
$ws.Range("A23").Value = "Chiết khấu thu nợ tại LONG XUYÊN"
$ws.Range("B23").Value = 0
$ws.Range("A24").Value = "Ứng lương tại LONG XUYÊN"
$ws.Range("B24").Value = -0
$ws.Range("A25").Value = "Tổng công tại SÓC TRĂNG"
$ws.Range("B25").Value = 0
$ws.Range("A26").Value = "Phụ cấp tại SÓC TRĂNG"
$ws.Range("B26").Value = 0
$ws.Range("A27").Value = "Lương công tác tại SÓC TRĂNG"
$ws.Range("B27").Value = 0
$ws.Range("A28").Value = "Lương cơ bản tại SÓC TRĂNG"
$ws.Range("B28").Value = 17250000
$ws.Range("A29").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws.Range("B29").Value = 0
$ws.Range("A30").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws.Range("B30").Value = 0
$ws.Range("A31").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws.Range("B31").Value = 0
$ws.Range("A32").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws.Range("B32").Value = 0
$ws.Range("A33").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws.Range("B34").Value = 0
$ws.Range("A35").Value = "Chiết khấu thu nợ tại SÓC TRĂNG"
$ws.Range("B35").Value = 0
$ws.Range("A36").Value = "Ứng lương tại SÓC TRĂNG"
$ws.Range("B36").Value = -0
$ws.Range("A37").Value = "Tổng lương tại CẦN THƠ"
$ws.Range("B37").Value = 15271000
$ws.Range("A38").Value = "Tổng lương tại LONG XUYÊN"
$ws.Range("B38").Value = 12410000
$ws.Range("A39").Value = "Tổng lương tại SÓC TRĂNG"
$ws.Range("B39").Value = 17250000
$ws.Range("A40").Value = "Tổng lương"
$ws.Range("B40").Value = 44931000
